$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -3
    3  = 0
    4  = -2
    5  = 0
    6  = -1
    7  = 1
    8  = -3
    9  = -3
    10 = -1
    11 = -5
    12 = 1
    14 = 2
    16 = 7
    17 = 1
    18 = -1
    19 = 2
    21 = 0
    22 = 5
    23 = -12
    24 = -6
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
